$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New crime data collected: weekly CompStat report rolls forward one week
# (Volume 31 Number 45 -> 46; week of 11/4-11/10/2024 -> 11/11-11/17/2024)
# and the Week-to-Date / 28-Day / YTD / 2-Year precinct figures are refreshed
# with the newly collected counts.

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/11/2024  Through  11/17/2024"

# --- Cells whose value TYPE flips between number and text this week ---
# (a precinct stat that was '0'/blank text now has a real count, or vice
# versa) -- copy number-format only from an already-correctly-styled cell
# in the same row, then write the new value so style+type both land right.
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "0"

$ws.Range("I22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1

$ws.Range("I22").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1

$ws.Range("D27").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("G27").Value = "0"

$ws.Range("E27").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("H27").Value = "***.*"

$ws.Range("G31").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 1

$ws.Range("H31").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100

# --- Remaining precinct figures: counts + percent changes refreshed ---
# --- Row 15: percentage figures update (M/N) ---
$ws.Range("M15").Value = 75
$ws.Range("N15").Value = 9.375

# --- Row 16 ---
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -21.739130434782
$ws.Range("I16").Value = 220
$ws.Range("J16").Value = 253
$ws.Range("K16").Value = -13.043478260869
$ws.Range("L16").Value = -7.563025210084
$ws.Range("M16").Value = -20.289855072463
$ws.Range("N16").Value = -64.052287581699

# --- Row 17 ---
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 12.5
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = -27.5
$ws.Range("I17").Value = 413
$ws.Range("J17").Value = 376
$ws.Range("K17").Value = 9.840425531914
$ws.Range("L17").Value = 26.299694189602
$ws.Range("M17").Value = 79.565217391304
$ws.Range("N17").Value = 34.527687296416

# --- Row 18 ---
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -8.333333333333
$ws.Range("I18").Value = 145
$ws.Range("J18").Value = 211
$ws.Range("K18").Value = -31.279620853080
$ws.Range("L18").Value = 27.192982456140
$ws.Range("M18").Value = -57.602339181286
$ws.Range("N18").Value = -88.948170731707

# --- Row 19 ---
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -30.769230769230
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = -20
$ws.Range("I19").Value = 750
$ws.Range("J19").Value = 601
$ws.Range("K19").Value = 24.792013311148
$ws.Range("L19").Value = 37.111517367458
$ws.Range("M19").Value = 87.969924812030
$ws.Range("N19").Value = 44.230769230769

# --- Row 20 ---
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = 8.333333333333
$ws.Range("F20").Value = 35
$ws.Range("G20").Value = 42
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 403
$ws.Range("J20").Value = 466
$ws.Range("K20").Value = -13.519313304721
$ws.Range("L20").Value = 27.936507936507
$ws.Range("M20").Value = 100.497512437811
$ws.Range("N20").Value = -76.238207547169

# --- Row 21 ---
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 2.631578947368
$ws.Range("F21").Value = 146
$ws.Range("G21").Value = 177
$ws.Range("H21").Value = -17.514124293785
$ws.Range("I21").Value = 1969
$ws.Range("J21").Value = 1933
$ws.Range("K21").Value = 1.862390067252
$ws.Range("L21").Value = 25.975687779910
$ws.Range("M21").Value = 33.220568335588
$ws.Range("N21").Value = -56.156757960365

# --- Row 22 ---
$ws.Range("I22").Value = 22
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 10
$ws.Range("M22").Value = 29.411764705882

# --- Row 23 ---
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -55.555555555555
$ws.Range("I23").Value = 95
$ws.Range("J23").Value = 107
$ws.Range("K23").Value = -11.214953271028
$ws.Range("L23").Value = -1.041666666666
$ws.Range("M23").Value = 66.666666666666

# --- Row 24 ---
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 4.545454545454
$ws.Range("F24").Value = 70
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = -33.962264150943
$ws.Range("I24").Value = 1107
$ws.Range("J24").Value = 1423
$ws.Range("K24").Value = -22.206605762473
$ws.Range("L24").Value = -7.980049875311
$ws.Range("M24").Value = 34.181818181818

# --- Row 25 ---
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 42.857142857142
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = -55.102040816326
$ws.Range("I25").Value = 423
$ws.Range("J25").Value = 612
$ws.Range("K25").Value = -30.882352941176
$ws.Range("L25").Value = -11.320754716981

# --- Row 26 ---
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -42.857142857142
$ws.Range("F26").Value = 49
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 11.363636363636
$ws.Range("I26").Value = 504
$ws.Range("J26").Value = 473
$ws.Range("K26").Value = 6.553911205074
$ws.Range("L26").Value = 7.922912205567
$ws.Range("M26").Value = -12.651646447140

# --- Row 27 ---
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 7
$ws.Range("I27").Value = 41
$ws.Range("K27").Value = 32.258064516129
$ws.Range("L27").Value = 20.588235294117

# --- Row 28 ---
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 62
$ws.Range("J28").Value = 62
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 26.530612244898

# --- Row 29 ---
$ws.Range("M29").Value = -32
$ws.Range("N29").Value = -62.222222222222

# --- Row 30 ---
$ws.Range("M30").Value = -28.571428571428
$ws.Range("N30").Value = -62.5

# --- Row 31 ---
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 4
$ws.Range("K31").Value = -75

# --- Row 33 ---
$ws.Range("L33").Value = -50
